$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-12 Monday" "2024-02-13 Tuesday"

Replace-Text "61×98=5978" "16×88=1408"
Replace-Text "86×14=1204" "48×79=3792"
Replace-Text "55×29=1595" "82×95=7790"
Replace-Text "60×85=5100" "50×42=2100"
Replace-Text "87×35=3045" "24×71=1704"

Replace-Text "20×61=1220" "70×96=6720"
Replace-Text "24×41=984" "36×31=1116"
Replace-Text "90×85=7650" "77×78=6006"
Replace-Text "63×93=5859" "41×40=1640"
Replace-Text "73×71=5183" "37×33=1221"

Replace-Text "61×65=3965" "15×12=180"
Replace-Text "41×93=3813" "16×78=1248"
Replace-Text "58×57=3306" "18×29=522"
Replace-Text "87×20=1740" "98×20=1960"
Replace-Text "39×18=702" "58×30=1740"

Replace-Text "12×40=480" "63×35=2205"
Replace-Text "85×26=2210" "21×46=966"
Replace-Text "62×26=1612" "87×80=6960"
Replace-Text "91×95=8645" "52×84=4368"
Replace-Text "94×57=5358" "84×66=5544"

Replace-Text "16×48=768" "95×79=7505"
Replace-Text "46×18=828" "30×48=1440"
Replace-Text "24×50=1200" "78×68=5304"
Replace-Text "82×12=984" "28×34=952"
Replace-Text "98×13=1274" "38×53=2014"
